$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header (row 1) cell values — comment.* labels for the new columns
$ws.Range("H1").Value = '<%=comment.create_usr_id_lbl%><%selectList.create_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.create_usr_id.join(",") }"` })%>'
$ws.Range("I1").Value = '<%=comment.create_time_lbl%>'
$ws.Range("J1").Value = '<%=comment.update_usr_id_lbl%><%selectList.update_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.update_usr_id.join(",") }"` })%>'
$ws.Range("K1").Value = '<%=comment.update_time_lbl%>'

# New data (row 2) cell values — model.* template tokens for the new columns
$ws.Range("H2").Value = '<%=model.create_usr_id_lbl%>'
$ws.Range("I2").Value = '<%~model.create_time ? new Date(model.create_time) : ""%>'
$ws.Range("J2").Value = '<%=model.update_usr_id_lbl%>'
$ws.Range("K2").Value = '<%~model.update_time ? new Date(model.update_time) : ""%>'
